$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "W2YVB8"
$ws.Range("B21").Value = "Caja de Mantenimiento T6714"
$ws.Range("C21").Value = "WF C814 C860 C869 C878 C879C8190 C8610 C8690"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 300000
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 2
$ws.Range("H21").Formula = "=(E21-D21)*G21"
$ws.Range("I21").Formula = "=D21*F21"
$ws.Range("J21").Value = 0
